# Updated cryptos list on Mon Jun 24 04:15:06 UTC 2024 with GitHub Actions
#
# Refreshes the Price (col D) and Volume(1h) (col E) figures for every coin
# row, plus the Litecoin/Dai row swap (rows 25-26: names, links, price and
# volume all exchange places since the ranking order changed).
#
# Values that look like a plain number (e.g. "1.00", "22.80", "0.0000173")
# are prefixed with a leading single-quote, PowerShell's own escape for a
# literal apostrophe inside a '...' string. That mirrors typing '1.00 into
# Excel by hand: the apostrophe forces text storage instead of being
# auto-coerced to a numeric value (which would silently drop the formatting,
# e.g. "22.80" -> 22.8 or "0.0000173" -> 1.73E-05). Values that are already
# unambiguous as text (contain a second "." like "62.740.72", or the "%"
# strings) are left unprefixed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col='D'; Value='62.740.72'},
    @{Row=2; Col='E'; Value='  -2.58%  '},
    @{Row=3; Col='D'; Value='3.392.46'},
    @{Row=3; Col='E'; Value='  -3.57%  '},
    @{Row=4; Col='E'; Value='  +0.06%  '},
    @{Row=5; Col='D'; Value='''574.84'},
    @{Row=5; Col='E'; Value='  -2.99%  '},
    @{Row=6; Col='D'; Value='''125.76'},
    @{Row=6; Col='E'; Value='  -6.73%  '},
    @{Row=7; Col='E'; Value='  +0.04%  '},
    @{Row=8; Col='D'; Value='3.390.05'},
    @{Row=8; Col='E'; Value='  -3.61%  '},
    @{Row=9; Col='D'; Value='''0.475'},
    @{Row=9; Col='E'; Value='  -2.82%  '},
    @{Row=10; Col='D'; Value='''7.37'},
    @{Row=10; Col='E'; Value='  -2.93%  '},
    @{Row=11; Col='D'; Value='''0.121'},
    @{Row=11; Col='E'; Value='  -2.49%  '},
    @{Row=12; Col='D'; Value='''0.379'},
    @{Row=12; Col='E'; Value='  -2.43%  '},
    @{Row=13; Col='D'; Value='3.981.08'},
    @{Row=13; Col='E'; Value='  -3.28%  '},
    @{Row=14; Col='E'; Value='  -0.94%  '},
    @{Row=15; Col='D'; Value='3.401.91'},
    @{Row=15; Col='E'; Value='  -3.32%  '},
    @{Row=16; Col='D'; Value='''0.0000173'},
    @{Row=16; Col='E'; Value='  -4.54%  '},
    @{Row=17; Col='D'; Value='62.763.46'},
    @{Row=17; Col='E'; Value='  -2.52%  '},
    @{Row=18; Col='D'; Value='''24.81'},
    @{Row=18; Col='E'; Value='  -4.21%  '},
    @{Row=19; Col='D'; Value='''9.52'},
    @{Row=19; Col='E'; Value='  -4.29%  '},
    @{Row=20; Col='D'; Value='''5.68'},
    @{Row=20; Col='E'; Value='  -1.55%  '},
    @{Row=21; Col='D'; Value='''13.19'},
    @{Row=21; Col='E'; Value='  -2.99%  '},
    @{Row=22; Col='D'; Value='''376.45'},
    @{Row=22; Col='E'; Value='  -4.65%  '},
    @{Row=23; Col='D'; Value='''0.558'},
    @{Row=23; Col='E'; Value='  -3.24%  '},
    @{Row=24; Col='D'; Value='3.530.52'},
    @{Row=24; Col='E'; Value='  -3.50%  '},
    @{Row=25; Col='B'; Value='Dai'},
    @{Row=25; Col='C'; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Row=25; Col='D'; Value='''1.00'},
    @{Row=25; Col='E'; Value='  -0.03%  '},
    @{Row=26; Col='B'; Value='Litecoin'},
    @{Row=26; Col='C'; Value='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{Row=26; Col='D'; Value='''72.39'},
    @{Row=26; Col='E'; Value='  -3.24%  '},
    @{Row=27; Col='D'; Value='''0.0000108'},
    @{Row=27; Col='E'; Value='  -7.80%  '},
    @{Row=28; Col='D'; Value='''0.999'},
    @{Row=28; Col='E'; Value='  -0.05%  '},
    @{Row=29; Col='E'; Value='  -5.67%  '},
    @{Row=30; Col='E'; Value='  -4.62%  '},
    @{Row=31; Col='D'; Value='''7.86'},
    @{Row=31; Col='E'; Value='  -5.51%  '},
    @{Row=32; Col='D'; Value='''1.40'},
    @{Row=32; Col='E'; Value='  -4.44%  '},
    @{Row=33; Col='D'; Value='''0.150'},
    @{Row=33; Col='E'; Value='  -5.29%  '},
    @{Row=34; Col='E'; Value='  -0.04%  '},
    @{Row=35; Col='D'; Value='3.420.93'},
    @{Row=35; Col='E'; Value='  -3.59%  '},
    @{Row=36; Col='D'; Value='''22.80'},
    @{Row=36; Col='E'; Value='  -2.84%  '},
    @{Row=37; Col='E'; Value='  -1.34%  '},
    @{Row=38; Col='D'; Value='''6.75'},
    @{Row=38; Col='E'; Value='  -3.20%  '},
    @{Row=39; Col='D'; Value='''164.50'},
    @{Row=39; Col='E'; Value='  -1.67%  '},
    @{Row=40; Col='D'; Value='''1.49'},
    @{Row=40; Col='E'; Value='  -4.52%  '},
    @{Row=41; Col='D'; Value='''0.0760'},
    @{Row=41; Col='E'; Value='  -4.07%  '},
    @{Row=42; Col='E'; Value='  +0.09%  '},
    @{Row=43; Col='D'; Value='''0.775'},
    @{Row=43; Col='E'; Value='  -4.57%  '},
    @{Row=44; Col='D'; Value='''41.44'},
    @{Row=44; Col='E'; Value='  -2.23%  '},
    @{Row=45; Col='D'; Value='''4.28'},
    @{Row=45; Col='E'; Value='  -4.01%  '},
    @{Row=46; Col='D'; Value='''1.57'},
    @{Row=46; Col='E'; Value='  -5.98%  '},
    @{Row=47; Col='D'; Value='''22.80'},
    @{Row=47; Col='E'; Value='  -10.92%  '},
    @{Row=48; Col='D'; Value='''1.08'},
    @{Row=48; Col='E'; Value='  -7.89%  '},
    @{Row=49; Col='D'; Value='''6.66'},
    @{Row=49; Col='E'; Value='  -2.23%  '},
    @{Row=50; Col='D'; Value='2.241.80'},
    @{Row=50; Col='E'; Value='  -7.13%  '},
    @{Row=51; Col='D'; Value='''0.855'},
    @{Row=51; Col='E'; Value='  -5.07%  '}
)

$colMap = @{ A=1; B=2; C=3; D=4; E=5 }

foreach ($u in $updates) {
    $colNum = $colMap[$u.Col]
    $ws.Cells.Item($u.Row, $colNum).Value = $u.Value
}

